# Apply updated cryptocurrency price/volume data (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns B/C/D/E are stored as plain text in the workbook (prices use
# "." as a thousands separator in some rows, e.g. "89.858.17"). Values that Excel
# would otherwise auto-convert to a number (e.g. "1.00", "3.71") are written with a
# leading apostrophe to force text, then the style is reset to Normal so no stray
# number-format is left applied to the cell.

$ws.Range('D2').Value = '89.858.17'
$ws.Range('E2').Value = '  -1.71%  '

$ws.Range('D3').Value = '3.058.62'
$ws.Range('E3').Value = '  -1.90%  '

$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').Value = "'215.41"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.07%  '

$ws.Range('D6').Value = "'609.60"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.92%  '

$ws.Range('E7').Value = '  +18.19%  '

$ws.Range('D8').Value = "'0.346"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -9.16%  '

$ws.Range('D9').Value = "'1.00"
$ws.Range('D9').Style = 'Normal'

$ws.Range('D10').Value = '3.056.89'
$ws.Range('E10').Value = '  -1.78%  '

$ws.Range('D11').Value = "'0.715"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.25%  '

$ws.Range('D12').Value = "'0.194"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.70%  '

$ws.Range('D13').Value = "'0.0000237"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -6.57%  '

$ws.Range('E14').Value = '  +0.64%  '

$ws.Range('D15').Value = "'33.91"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.33%  '

$ws.Range('D16').Value = '89.664.14'
$ws.Range('E16').Value = '  -1.60%  '

$ws.Range('D17').Value = '3.626.19'
$ws.Range('E17').Value = '  -1.75%  '

$ws.Range('D18').Value = '3.053.02'
$ws.Range('E18').Value = '  -2.12%  '

$ws.Range('D19').Value = "'3.71"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.68%  '

$ws.Range('D20').Value = "'14.05"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.23%  '

$ws.Range('D21').Value = "'444.92"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.91%  '

$ws.Range('D22').Value = "'0.0000200"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -13.23%  '

$ws.Range('D23').Value = "'5.37"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.09%  '

$ws.Range('D24').Value = "'8.76"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.23%  '

$ws.Range('D25').Value = "'5.76"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.17%  '

$ws.Range('D26').Value = "'90.42"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.00%  '

$ws.Range('D27').Value = "'11.76"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.86%  '

$ws.Range('D28').Value = '3.230.17'
$ws.Range('E28').Value = '  -1.34%  '

$ws.Range('E29').Value = '  -0.04%  '

$ws.Range('D30').Value = "'9.08"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.20%  '

$ws.Range('D31').Value = "'0.158"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.19%  '

$ws.Range('B32').Value = 'Binance-PegBSC-USD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D32').Value = "'1.00"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +9.95%  '

$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = "'27.98"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +20.27%  '

$ws.Range('B34').Value = 'Stellar'
$ws.Range('C34').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D34').Value = "'0.200"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +40.39%  '

$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = "'0.145"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.75%  '

$ws.Range('B36').Value = 'Bittensor'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D36').Value = "'492.50"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.95%  '

$ws.Range('D37').Value = "'1.87"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.92%  '

$ws.Range('D38').Value = "'6.70"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.63%  '

$ws.Range('D39').Value = "'1.27"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.76%  '

$ws.Range('D40').Value = "'3.38"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -12.30%  '

$ws.Range('D41').Value = "'0.421"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +11.99%  '

$ws.Range('D42').Value = "'22.18"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.53%  '

$ws.Range('E43').Value = '  -0.02%  '

$ws.Range('D44').Value = "'0.0842"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +10.21%  '

$ws.Range('D45').Value = "'1.91"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.14%  '

$ws.Range('D46').Value = "'3.02"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +15.17%  '

$ws.Range('D47').Value = "'147.39"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.67%  '

$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').Value = "'0.677"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +10.32%  '

$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').Value = "'4.51"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.77%  '

$ws.Range('D50').Value = "'44.61"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.30%  '

$ws.Range('B51').Value = 'ImmutableX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D51').Value = "'1.30"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.15%  '
